# save data done + era data updated
# Add a new "Save" column (H) to the sheet, mirroring the header style
# used by the existing columns, and populate the data row with the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it reuses the same bold/border/centered style, then
# set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column in row 2.
$ws.Range("H2").Value = 1
